$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New audit-trail header columns (F:K) ---
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# --- Fix native_name for French language row (re-import artifact) ---
$ws.Range("D2").Value = "franÃ§ais"

# --- is_active is now a real boolean TRUE instead of text "TRUE" ---
$ws.Range("E2").Value = $true

# --- New audit-trail values for the French row (F2:K2) ---
$ws.Range("F2").Value = "superadmin"
$ws.Range("G2").Value = 45079.577112638886
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = "NULL"

# --- Restore the selected cell left by the authoring session ---
$ws.Range("E7").Select()
